# Update the "想去人数" (interest count) column F on the "展览" and
# "全部类型" worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 104
$ws1.Range("F10").Value = 164
$ws1.Range("F14").Value = 292
$ws1.Range("F16").Value = 3837
$ws1.Range("F24").Value = 2901
$ws1.Range("F32").Value = 2096
$ws1.Range("F33").Value = 942
$ws1.Range("F36").Value = 529
$ws1.Range("F37").Value = 319
$ws1.Range("F41").Value = 998
$ws1.Range("F46").Value = 323
$ws1.Range("F47").Value = 259

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 104
$ws4.Range("F9").Value  = 164
$ws4.Range("F13").Value = 292
$ws4.Range("F15").Value = 3837
$ws4.Range("F23").Value = 2901
$ws4.Range("F34").Value = 2096
$ws4.Range("F36").Value = 942
$ws4.Range("F38").Value = 529
$ws4.Range("F39").Value = 319
$ws4.Range("F41").Value = 998
$ws4.Range("F46").Value = 323
$ws4.Range("F48").Value = 259
